$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register the size-8 "Aptos Narrow" font that the phonetic-info settings
# (Format > Phonetic Information) reference as fontId="1". We stamp it onto a
# scratch cell outside the used range and then clear that cell completely so
# the new font is added to the shared styles table without leaving any cell
# pointing at it (and without perturbing the sheet's used range/dimension).
$scratch = $ws.Range("Z1")
$scratch.Font.Size = 8
$scratch.Clear()

# Append the nine new leaderboard rows (test1..test9), each starting at 0.
$names = @("test1", "test2", "test3", "test4", "test5", "test6", "test7", "test8", "test9")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 5 + $i
    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("B$row").Value = 0
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
}

# Match the saved selection state (active cell just below the new data).
$ws.Range("B14").Select()
